$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (rows 1-8, columns A-C) before rewriting.
$ws.Range("A1:C8").ClearContents()

# Write new header row
$ws.Range("A1").Value = "Picture1"
$ws.Range("B1").Value = "Picture2"
$ws.Range("C1").Value = "Locations"
$ws.Range("D1").Value = "ProbeLocation"

# Row 2
$ws.Range("A2").Value = "smile.jpg"
$ws.Range("B2").Value = "frown.jpg"
$ws.Range("C2").Value = "[0.35, 0]"
$ws.Range("D2").Value = "smile"

# Row 3
$ws.Range("A3").Value = "frown.jpg"
$ws.Range("B3").Value = "smile.jpg"
$ws.Range("C3").Value = "[-0.35, 0]"
$ws.Range("D3").Value = "frown"

# Update selection to reflect the new active cell
$ws.Range("A2").Select()
